$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (date rolled forward one day)
$ws.Name = "Through 2022-07-23"

# Update the "July" label to reflect new as-of date
$ws.Range("A8").Value = "July (through 07-23)"

# Update July row (row 8) values
$ws.Range("B8").Value = 27
$ws.Range("D8").Value = 51
$ws.Range("E8").Value = 58
$ws.Range("F8").Value = 36
$ws.Range("G8").Value = 98
$ws.Range("H8").Value = 114
$ws.Range("I8").Value = 133

# Update Total row (row 9) values
$ws.Range("B9").Value = 152
$ws.Range("D9").Value = 441
$ws.Range("E9").Value = 411
$ws.Range("F9").Value = 287
$ws.Range("G9").Value = 570
$ws.Range("H9").Value = 874
$ws.Range("I9").Value = 939
